$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.449.88'
$ws.Range('E2').Value = '  +1.79%  '

$ws.Range('D3').Value = '2.423.87'
$ws.Range('E3').Value = '  -0.72%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '319.83'
$ws.Range('E5').Value = '  +3.43%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.94'
$ws.Range('E6').Value = '  +1.60%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.516'
$ws.Range('E7').Value = '  +0.80%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.530'
$ws.Range('E9').Value = '  +4.58%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.74'
$ws.Range('E10').Value = '  +0.88%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0801'
$ws.Range('E11').Value = '  +0.06%  '

$ws.Range('E12').Value = '  -1.31%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.20'
$ws.Range('E13').Value = '  -2.95%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.98'
$ws.Range('E14').Value = '  +0.11%  '

$ws.Range('D15').Value = '2.803.55'
$ws.Range('E15').Value = '  -0.60%  '

$ws.Range('D16').Value = '2.431.18'
$ws.Range('E16').Value = '  -0.67%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.836'
$ws.Range('E17').Value = '  -0.55%  '

$ws.Range('D18').Value = '45.349.95'
$ws.Range('E18').Value = '  +1.79%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.32'
$ws.Range('E19').Value = '  -1.09%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.38'
$ws.Range('E20').Value = '  -0.65%  '

$ws.Range('D21').Value = '0.0₃0927'
$ws.Range('E21').Value = '  +2.00%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.41'
$ws.Range('E22').Value = '  +3.70%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.37'
$ws.Range('E23').Value = '  +0.68%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '245.63'
$ws.Range('E24').Value = '  +1.84%  '

$ws.Range('E25').Value = '  -0.16%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.63'
$ws.Range('E27').Value = '  +1.42%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -2.61%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.65'
$ws.Range('E29').Value = '  -0.24%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '33.36'
$ws.Range('E30').Value = '  +0.40%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '49.24'
$ws.Range('E31').Value = '  +0.92%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  +6.34%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.17'
$ws.Range('E33').Value = '  +3.29%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.25'
$ws.Range('E34').Value = '  +0.75%  '

$ws.Range('E35').Value = '  +0.15%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0757'
$ws.Range('E36').Value = '  -1.31%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.51'
$ws.Range('E37').Value = '  -0.17%  '

$ws.Range('E38').Value = '  -1.84%  '

$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.91'
$ws.Range('E39').Value = '  -0.17%  '

$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '127.94'
$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('E41').Value = '  -3.20%  '

$ws.Range('E42').Value = '  +1.25%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.70'
$ws.Range('E43').Value = '  -5.09%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0291'
$ws.Range('E44').Value = '  +0.19%  '

$ws.Range('D45').Value = '1.951.68'
$ws.Range('E45').Value = '  -0.19%  '

$ws.Range('E46').Value = '  -2.80%  '

$ws.Range('E47').Value = '  +0.61%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.81'
$ws.Range('E48').Value = '  +8.03%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.12'
$ws.Range('E49').Value = '  -5.07%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '76.89'
$ws.Range('E50').Value = '  +3.96%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.81'
$ws.Range('E51').Value = '  +3.38%  '
